# CORREÇÃO DO BUG QUANDO O NUMERO É MAIOR QUUE O TABULEIRO
#
# 1) Remove the blank default "Sheet" worksheet (it was never used).
# 2) Update "configurações" values: Linhas=3, Colunas=3, Dificuldade=1.
# 3) Replace the "jogo" board (previously a random 5x5 grid) with a new
#    3x3 grid of values, matching the smaller board size above.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1) Delete the blank "Sheet" worksheet.
[void]$wb.Worksheets.Item("Sheet").Delete()

# 2) Update the "configurações" sheet values.
$cfg = $wb.Worksheets.Item("configurações")
$cfg.Range("A1").Value = "Linhas"
$cfg.Range("B1").Value = 3
$cfg.Range("A2").Value = "Colunas"
$cfg.Range("B2").Value = 3
$cfg.Range("A3").Value = "Dificuldade"
$cfg.Range("B3").Value = 1

# 3) Replace the "jogo" sheet board: shrink from a 5x5 grid to a new 3x3 grid.
$jogo = $wb.Worksheets.Item("jogo")
$jogo.Cells.Clear()

$values = @(
    @("0", "0", "0"),
    @("1", "1", "0"),
    @("-1", "1", "0")
)

# Write each value through a text-literal formula first (="0") so Excel
# commits it as text rather than re-interpreting the numeric-looking
# string as a number, then collapse the formulas down to plain values via
# copy / paste-values (xlPasteValues = -4163).
for ($r = 0; $r -lt 3; $r++) {
    for ($c = 0; $c -lt 3; $c++) {
        $jogo.Cells.Item($r + 1, $c + 1).Formula = '="' + $values[$r][$c] + '"'
    }
}

$boardRange = $jogo.Range("A1:C3")
[void]$boardRange.Copy()
[void]$boardRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false

$excel.DisplayAlerts = $true
